$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.4593233333333334
$ws.Range("H2").Value = 1.37797
$ws.Range("I2").Value = 0.015538272766109
$ws.Range("J2").Value = 0.015538272766109
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 110.8604276666667
$ws.Range("N2").Value = 332.581283
$ws.Range("O2").Value = 0.2509786052589675
$ws.Range("P2").Value = 0.2509786052589675
$ws.Range("Q2").Value = 50.92078117061222
$ws.Range("R2").Value = 458.28703053551
$ws.Range("S2").Value = 0.003899774026971436
$ws.Range("T2").Value = 0.003899774026971437

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.4593233333333334
$ws.Range("H3").Value = 1.37797
$ws.Range("I3").Value = 0.015538272766109
$ws.Range("J3").Value = 0.015538272766109
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 184.841802
$ws.Range("N3").Value = 554.525406
$ws.Range("O3").Value = 0.4184661617850055
$ws.Range("P3").Value = 0.4184661617850055
$ws.Range("Q3").Value = 84.90215263398001
$ws.Range("R3").Value = 764.1193737058201
$ws.Range("S3").Value = 0.006502241365202114
$ws.Range("T3").Value = 0.006502241365202115

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.4593233333333334
$ws.Range("H4").Value = 1.37797
$ws.Range("I4").Value = 0.015538272766109
$ws.Range("J4").Value = 0.015538272766109
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 95.23175666666667
$ws.Range("N4").Value = 285.69527
$ws.Range("O4").Value = 0.2155966197102082
$ws.Range("P4").Value = 0.2155966197102082
$ws.Range("Q4").Value = 43.74216791132223
$ws.Range("R4").Value = 393.6795112019
$ws.Range("S4").Value = 0.003349999084508288
$ws.Range("T4").Value = 0.003349999084508289

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.4593233333333334
$ws.Range("H5").Value = 1.37797
$ws.Range("I5").Value = 0.015538272766109
$ws.Range("J5").Value = 0.015538272766109
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 50.778675
$ws.Range("N5").Value = 152.336025
$ws.Range("O5").Value = 0.1149586132458188
$ws.Range("P5").Value = 0.1149586132458188
$ws.Range("Q5").Value = 23.32383026325
$ws.Range("R5").Value = 209.91447236925
$ws.Range("S5").Value = 0.001786258289427164
$ws.Range("T5").Value = 0.001786258289427164

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.94315233333334
$ws.Range("H6").Value = 35.829457
$ws.Range("I6").Value = 0.4040203167903319
$ws.Range("J6").Value = 0.4040203167903319
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 110.8604276666667
$ws.Range("N6").Value = 332.581283
$ws.Range("O6").Value = 0.2509786052589675
$ws.Range("P6").Value = 0.2509786052589675
$ws.Range("Q6").Value = 1324.022975361481
$ws.Range("R6").Value = 11916.20677825333
$ws.Range("S6").Value = 0.1014004556043237
$ws.Range("T6").Value = 0.1014004556043237

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.94315233333334
$ws.Range("H7").Value = 35.829457
$ws.Range("I7").Value = 0.4040203167903319
$ws.Range("J7").Value = 0.4040203167903319
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 184.841802
$ws.Range("N7").Value = 554.525406
$ws.Range("O7").Value = 0.4184661617850055
$ws.Range("P7").Value = 0.4184661617850055
$ws.Range("Q7").Value = 2207.593798853839
$ws.Range("R7").Value = 19868.34418968454
$ws.Range("S7").Value = 0.1690688312504122
$ws.Range("T7").Value = 0.1690688312504122

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.94315233333334
$ws.Range("H8").Value = 35.829457
$ws.Range("I8").Value = 0.4040203167903319
$ws.Range("J8").Value = 0.4040203167903319
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 95.23175666666667
$ws.Range("N8").Value = 285.69527
$ws.Range("O8").Value = 0.2155966197102082
$ws.Range("P8").Value = 0.2155966197102082
$ws.Range("Q8").Value = 1137.367376840932
$ws.Range("R8").Value = 10236.30639156839
$ws.Range("S8").Value = 0.08710541459424304
$ws.Range("T8").Value = 0.08710541459424305

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.94315233333334
$ws.Range("H9").Value = 35.829457
$ws.Range("I9").Value = 0.4040203167903319
$ws.Range("J9").Value = 0.4040203167903319
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 50.778675
$ws.Range("N9").Value = 152.336025
$ws.Range("O9").Value = 0.1149586132458188
$ws.Range("P9").Value = 0.1149586132458188
$ws.Range("Q9").Value = 606.4574508098251
$ws.Range("R9").Value = 5458.117057288426
$ws.Range("S9").Value = 0.04644561534135295
$ws.Range("T9").Value = 0.04644561534135296

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.844759
$ws.Range("H10").Value = 38.534277
$ws.Range("I10").Value = 0.4345204227020912
$ws.Range("J10").Value = 0.4345204227020912
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 110.8604276666667
$ws.Range("N10").Value = 332.581283
$ws.Range("O10").Value = 0.2509786052589675
$ws.Range("P10").Value = 0.2509786052589675
$ws.Range("Q10").Value = 1423.975476015266
$ws.Range("R10").Value = 12815.77928413739
$ws.Range("S10").Value = 0.1090553296463078
$ws.Range("T10").Value = 0.1090553296463078

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.844759
$ws.Range("H11").Value = 38.534277
$ws.Range("I11").Value = 0.4345204227020912
$ws.Range("J11").Value = 0.4345204227020912
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 184.841802
$ws.Range("N11").Value = 554.525406
$ws.Range("O11").Value = 0.4184661617850055
$ws.Range("P11").Value = 0.4184661617850055
$ws.Range("Q11").Value = 2374.248399815718
$ws.Range("R11").Value = 21368.23559834146
$ws.Range("S11").Value = 0.1818320935053422
$ws.Range("T11").Value = 0.1818320935053422

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 12.844759
$ws.Range("H12").Value = 38.534277
$ws.Range("I12").Value = 0.4345204227020912
$ws.Range("J12").Value = 0.4345204227020912
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 95.23175666666667
$ws.Range("N12").Value = 285.69527
$ws.Range("O12").Value = 0.2155966197102082
$ws.Range("P12").Value = 0.2155966197102082
$ws.Range("Q12").Value = 1223.228963529977
$ws.Range("R12").Value = 11009.06067176979
$ws.Range("S12").Value = 0.09368113432962168
$ws.Range("T12").Value = 0.09368113432962169

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 12.844759
$ws.Range("H13").Value = 38.534277
$ws.Range("I13").Value = 0.4345204227020912
$ws.Range("J13").Value = 0.4345204227020912
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 50.778675
$ws.Range("N13").Value = 152.336025
$ws.Range("O13").Value = 0.1149586132458188
$ws.Range("P13").Value = 0.1149586132458188
$ws.Range("Q13").Value = 652.2398427143251
$ws.Range("R13").Value = 5870.158584428926
$ws.Range("S13").Value = 0.0499518652208194
$ws.Range("T13").Value = 0.0499518652208194

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.313537
$ws.Range("H14").Value = 12.940611
$ws.Range("I14").Value = 0.1459209877414679
$ws.Range("J14").Value = 0.145920987741468
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 110.8604276666667
$ws.Range("N14").Value = 332.581283
$ws.Range("O14").Value = 0.2509786052589675
$ws.Range("P14").Value = 0.2509786052589675
$ws.Range("Q14").Value = 478.2005565759903
$ws.Range("R14").Value = 4303.805009183913
$ws.Range("S14").Value = 0.03662304598136452
$ws.Range("T14").Value = 0.03662304598136452

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.313537
$ws.Range("H15").Value = 12.940611
$ws.Range("I15").Value = 0.1459209877414679
$ws.Range("J15").Value = 0.145920987741468
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 184.841802
$ws.Range("N15").Value = 554.525406
$ws.Range("O15").Value = 0.4184661617850055
$ws.Range("P15").Value = 0.4184661617850055
$ws.Range("Q15").Value = 797.3219520736741
$ws.Range("R15").Value = 7175.897568663066
$ws.Range("S15").Value = 0.06106299566404892
$ws.Range("T15").Value = 0.06106299566404894

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.313537
$ws.Range("H16").Value = 12.940611
$ws.Range("I16").Value = 0.1459209877414679
$ws.Range("J16").Value = 0.145920987741468
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 95.23175666666667
$ws.Range("N16").Value = 285.69527
$ws.Range("O16").Value = 0.2155966197102082
$ws.Range("P16").Value = 0.2155966197102082
$ws.Range("Q16").Value = 410.7857059566634
$ws.Range("R16").Value = 3697.07135360997
$ws.Range("S16").Value = 0.03146007170183522
$ws.Range("T16").Value = 0.03146007170183523

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.313537
$ws.Range("H17").Value = 12.940611
$ws.Range("I17").Value = 0.1459209877414679
$ws.Range("J17").Value = 0.1459209877414679
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 50.778675
$ws.Range("N17").Value = 152.336025
$ws.Range("O17").Value = 0.1149586132458188
$ws.Range("P17").Value = 0.1149586132458188
$ws.Range("Q17").Value = 219.035693423475
$ws.Range("R17").Value = 1971.321240811275
$ws.Range("S17").Value = 0.01677487439421928
$ws.Range("T17").Value = 0.01677487439421928
